$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6668
$ws1.Range("F3").Value = 11
$ws1.Range("F4").Value = 419
$ws1.Range("F9").Value = 91
$ws1.Range("F13").Value = 397
$ws1.Range("F14").Value = 1291
$ws1.Range("F15").Value = 13
$ws1.Range("F16").Value = 3324
$ws1.Range("F17").Value = 17
$ws1.Range("F18").Value = 215
$ws1.Range("F19").Value = 1977
$ws1.Range("F20").Value = 91

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6668
$ws4.Range("F3").Value = 11
$ws4.Range("F4").Value = 419
$ws4.Range("F10").Value = 91
$ws4.Range("F14").Value = 397
$ws4.Range("F15").Value = 1291
$ws4.Range("F16").Value = 13
$ws4.Range("F17").Value = 3324
$ws4.Range("F18").Value = 17
$ws4.Range("F19").Value = 215
$ws4.Range("F20").Value = 1977
$ws4.Range("F21").Value = 91
